# Add columns I (I0) and J (IF) to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-20
$values = @(
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(6, 6),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(5, 6),
    @(5, 5),
    @(7, 7),
    @(6, 6),
    @(7, 7)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
